# adicionando matriz de rigidez e tentativa de jacobi e gauss-seidel
#
# Updates the element incidence table on the "Incidencia" sheet: several
# rows had their start/end node numbers (columns A and B) corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Incidencia")
$ws.Activate()

# row -> (A value, B value)
$updates = @{
    3  = @(1, 3)
    4  = @(2, 3)
    9  = @(4, 6)
    10 = @(5, 6)
    11 = @(4, 7)
    12 = @(6, 7)
    13 = @(6, 8)
    14 = @(7, 8)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
}

$ws.Range("B6").Select()
